$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the "Qui" (owner) column (D) for a handful of rows.
# New shared strings must be created in the same order Excel assigned them
# ("Emma" then "Flo") so the shared-strings table matches the target file.
$ws.Range("D15").Value = "Emma"
$ws.Range("D8").Value  = "Flo"
$ws.Range("D9").Value  = "Flo"
$ws.Range("D7").Value  = "YES"
$ws.Range("D11").Value = "YES"
$ws.Range("D12").Value = "YES"
